# pax02.xlsx — ongoing work, csv added, better views
#
# Rename three header cells on Sheet1 (the raw serial-number / date-format
# headers get tidied up), then move the active selection from the last
# data cell (J4) up to the first renamed header cell (F1).
#
# NOTE on ordering: Excel's shared-string table appends newly-introduced
# strings in the order the cells are actually written (and drops strings
# that become unused). To reproduce the exact shared-string order from the
# target workbook, the date headers (H1, I1) must be written before the
# serial-number header (F1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "ExpDate_Excel"
$ws.Range("I1").Value = "DoB_Excel"
$ws.Range("F1").Value = "serialNrID"

$ws.Range("F1").Select() | Out-Null
